# "added text to thesis" -- extend the 95th-percentile boxplot table with a
# new "Top whisker 95% (t)" row-sum column (G), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets its own standalone SUM formula.
$ws.Range("G2").Formula = "=SUM(B2:E2)"

# Rows 5-9 are entered as one fill-down, which Excel records as a shared
# formula group (anchored at G5, members G6:G9).
$ws.Range("G5:G9").Formula = "=SUM(B5:E5)"

# The interior helper rows (6,7,8) are cleared back out, leaving only the
# totals for the "Box mid line 50%" (row 5) and "max (t)" (row 9) rows - the
# cells stay in the sheet (carrying the number format) but without content.
$ws.Range("G6:G8").ClearContents()

# Give every new G cell in the table the same numeric style (format 0.000)
# used by the rest of the B:E data columns.
$ws.Range("G2:G9").NumberFormat = $ws.Range("E2").NumberFormat

# Reflect the new used range and the author's final on-screen selection.
$ws.Range("D40:F43").Select()
